# Correción archivo densidad poblacional
# Replace the "Population" column (B) values - which were stored as text
# strings with thousands separators (e.g. "16,787,941") or, for a few rows,
# as mis-scaled decimal numbers (e.g. 243.247 meant to be 243247) - with
# proper whole numbers, formatted with a plain integer number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$population = @(
    16787941,
    1055450,
    1247953,
    243247,
    64473,
    104099452,
    91276115,
    33406061,
    199812341,
    343709,
    25351462,
    72147030,
    27743338,
    32988134,
    31205576,
    1458545,
    112374333,
    3673917,
    61095297,
    84580777,
    60439692,
    35193978,
    41974218,
    72626809,
    68548437,
    25545198,
    10086292,
    2966889,
    2855794,
    6864602,
    1978502,
    610577,
    12541302,
    1097206,
    380581,
    1383727
)

for ($i = 0; $i -lt $population.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $population[$i]
}

$ws.Range("B2:B37").NumberFormat = "0"

# Restore portrait page setup and move the active selection back to E3
$ws.PageSetup.Orientation = 1
[void]$ws.Range("E3").Select()
